$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-4.377829468602984,-3.038932295093685,0.2056286196336669,-0.1561269739048682,0.3138575077444884),
    @(-4.135240465056005,-2.925761869611146,0.1526961756983606,-0.007081435956113867,0.1996924965362855),
    @(-3.91779440813451,-2.786980976820223,0.2116711397399643,0.06285134145811457,0.1409986714813057),
    @(-3.727207631739984,-2.63896785544131,0.266773802289977,0.1151279292582288,0.09726210584603523),
    @(-3.560370163304749,-2.486596470533827,0.3083680167545987,0.156400625316125,0.06290010432766577),
    @(-3.4140320213363,-2.333043362829101,0.3373751349543979,0.1880986218500402,0.03601604595621438),
    @(-3.285267012599882,-2.180799115878838,0.3556584128154492,0.2112694803727222,0.01526419832980576),
    @(-3.171504433694817,-2.031820051881835,0.365023877449405,0.2269773108623601,-0.0004787394824077064),
    @(-3.070510228545634,-1.887602301144848,0.3671011042269253,0.2362818665803226,-0.01216270276063274),
    @(-2.980358720832418,-1.749249412619633,0.3633261933465243,0.2401960799948154,-0.02058573907449598),
    @(-2.899401465639657,-1.617536338918757,0.3549447287318384,0.2396572970327526,-0.02641257492002178),
    @(-2.826235676702645,-1.492967745081127,0.3430228639462004,0.2355116619108308,-0.03019228911398673),
    @(-2.75967371208912,-1.375829432128243,0.3284626296547462,0.2285081468464749,-0.03237490830192973),
    @(-2.698714543450587,-1.266232600289749,0.3120191390451815,0.2192991240412868,-0.0333266207866517),
    @(-2.642517725653839,-1.16415126757514,0.2943181525803341,0.2084450816413542,-0.03334340822790796),
    @(-2.59038008719811,-1.069453473394109,0.2758730185441844,0.1964217272548262,-0.03266301278600346),
    @(-2.541715157163297,-0.9819270265325778,0.2571004096770639,0.1836282471547657,-0.03147525679828501),
    @(-2.496035212057282,-0.9013005735026718,0.2383345612355345,0.1703958948165927,-0.02993080535515241),
    @(-2.452935747855727,-0.827260716904461,0.219839908663442,0.1569963830630042,-0.02814850927639157),
    @(-2.412082143414136,-0.7594658365778045,0.2018221467515104,0.1436497688914796,-0.02622149068009805),
    @(-2.373198268962183,-0.6975571791258875,0.1844378057736476,0.1305316681338991,-0.02422214067678939),
    @(-2.336056798143359,-0.6411676954517767,0.1678024787530985,0.1177797353006681,-0.02220619385725165),
    @(-2.300470997127772,-0.589929027279597,0.1519978492892873,0.105499406126406,-0.02021603166578994),
    @(-2.266287784834241,-0.5434769747808197,0.1370776697934964,0.09376893726324381,-0.01828335001598595),
    @(-2.233381880963938,-0.5014557189097788,0.1230728316179136,0.08264379738313567,-0.01643130821148723),
    @(-2.201650881263391,-0.4635210232883112,0.1099956555985796,0.07216047257004851,-0.01467625813462307),
    @(-2.17101112093779,-0.4293426003682588,0.09784351672627689,0.06233975049621875,-0.01302913583962656),
    @(-2.141394206750376,-0.3986057938719949,0.08660190173614611,0.05318954536745196,-0.01149658270986926),
    @(-2.112744115783735,-0.3710127029166113,0.07624698436831955,0.04470732091579457,-0.01008185043184822),
    @(-2.085014774065904,-0.3462828516364005,0.06674779040933669,0.03688216305387852,-0.008785533194832459),
    @(-2.058168041386568,-0.3241534905576728,0.05806801357067486,0.02969654794681034,-0.007606161598544746),
    @(-2.032172039833269,-0.3043796016399409,0.05016753379413031,0.02312784564417318,-0.006540685525120827),
    @(-2.006999773084602,-0.2867336671206656,0.04300368158159788,0.01714959427347313,-0.00558486746229182),
    @(-1.982627991532646,-0.2710052525625544,0.03653228526271295,0.01173257521948482,-0.004733603210405914),
    @(-1.959036265084278,-0.2570004463902629,0.03070853255398958,0.006845715712852657,-0.003981183342925771),
    @(-1.936206231204822,-0.2445411913953329,0.02548767314745623,0.002456841792197529,-0.003321506021298277),
    @(-1.914120990588503,-0.2334645379349274,0.02082558523122898,-0.001466698368816741,-0.002748249622872337),
    @(-1.892764626914568,-0.2236218436582191,0.01667922564586455,-0.004957523327726419,-0.002255011986036276),
    @(-1.872121830599931,-0.2148779404147087,0.01300698069802322,-0.008047804180073876,-0.001835421797695067),
    @(-1.852177609392017,-0.2071102854145721,0.009768932387775904,-0.01076893393933568,-0.001483226655444758)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

$wb.Save()